# Update Price (D) and Volume(1h) (E) columns with refreshed crypto market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.507.38"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "'1.869.27"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'312.36"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "'0.4779"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "'0.3775"
$ws.Range("E8").Value = "  +2.93%  "
$ws.Range("D9").Value = "'0.07350"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "'0.9373"
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("D11").Value = "'20.73"
$ws.Range("E11").Value = "  +5.11%  "
$ws.Range("D12").Value = "'0.07859"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("D13").Value = "'1.856.22"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'5.437"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").Value = "'6.574"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("D16").Value = "'90.70"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D18").Value = "'0.000008903"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("D21").Value = "'27.520.02"
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("D22").Value = "'5.131"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").Value = "'10.72"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").Value = "'1.956"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "'153.79"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").Value = "'18.51"
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("D27").Value = "'2.026"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("D28").Value = "'115.89"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").Value = "'0.08927"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "'3.339"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").Value = "'1.215"
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").Value = "'4.614"
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("D34").Value = "'0.7534"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").Value = "'2.709"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'0.02051"
$ws.Range("E36").Value = "  +4.92%  "
$ws.Range("D37").Value = "'1.119"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").Value = "'3.003"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "'0.05275"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").Value = "'0.5343"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("D41").Value = "'7.087"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").Value = "'0.1525"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").Value = "'8.482"
$ws.Range("E43").Value = "  +3.35%  "
$ws.Range("D44").Value = "'10.68"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").Value = "'0.4807"
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("D46").Value = "'1.013"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Value = "'1.659"
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("D48").Value = "'102.87"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").Value = "'67.49"
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("D50").Value = "'0.06087"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").Value = "'0.9243"
$ws.Range("E51").Value = "  +4.45%  "
